# Updated cryptos list on Sat Aug 19 11:28:11 UTC 2023 with GitHub Actions
# Refresh price / 1h-volume-change figures scraped from coinranking.com,
# and re-rank Maker above VeChain to reflect the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as text (values may use "." as a thousands
# separator, e.g. "26.048.54"), so force text format before writing
# to avoid Excel auto-converting them to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.048.54"
$ws.Range("E2").Value = "  -2.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.667.92"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.66"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2655"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06391"
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.77"
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07446"
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.669.59"
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.511"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5825"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008548"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.28"
$ws.Range("E16").Value = "  -1.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.112.26"
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.938"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.41"
$ws.Range("E21").Value = "  +2.58%  "
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.76"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.611"
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("E26").Value = "  +2.19%  "
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06582"
$ws.Range("E28").Value = "  +14.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.340"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.511"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6117"
$ws.Range("E35").Value = "  +2.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.370"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.688"
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("E38").Value = "  +7.84%  "

# Rows 39 & 40 swap places: Maker moves to rank 39, VeChain to rank 40,
# each with refreshed price/volume figures.
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.091.87"
$ws.Range("E39").Value = "  +0.07%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01599"
$ws.Range("E40").Value = "  -1.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8713"
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.010"
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.14"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.816.65"
$ws.Range("E44").Value = "  -1.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000111"
$ws.Range("E45").Value = "  -5.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.34"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.078"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4287"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.036"
$ws.Range("E51").Value = "  +4.27%  "